$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 103, pushing the existing rows 103:126 down to 105:128.
$ws.Rows.Item(103).Resize(2).Insert()

# New row 103: Black Cat / Especial
$ws.Cells.Item(103, 1).Value = 9
$ws.Cells.Item(103, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(103, 3).Value = "Metropolitana"
$ws.Cells.Item(103, 4).Value = 44644
$ws.Cells.Item(103, 5).Value = 13
$ws.Cells.Item(103, 6).Value = "Fruta"
$ws.Cells.Item(103, 7).Value = 100103
$ws.Cells.Item(103, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(103, 9).Value = 100103002
$ws.Cells.Item(103, 10).Value = "Ciruela"
$ws.Cells.Item(103, 11).Value = "Black Cat"
$ws.Cells.Item(103, 12).Value = "Especial"
$ws.Cells.Item(103, 13).Value = 250
$ws.Cells.Item(103, 14).Value = 9000
$ws.Cells.Item(103, 15).Value = 9000
$ws.Cells.Item(103, 16).Value = 9000
$ws.Cells.Item(103, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(103, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(103, 19).Value = 600
$ws.Cells.Item(103, 20).Value = 15

# New row 104: Black Cat / Primera
$ws.Cells.Item(104, 1).Value = 9
$ws.Cells.Item(104, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(104, 3).Value = "Metropolitana"
$ws.Cells.Item(104, 4).Value = 44644
$ws.Cells.Item(104, 5).Value = 13
$ws.Cells.Item(104, 6).Value = "Fruta"
$ws.Cells.Item(104, 7).Value = 100103
$ws.Cells.Item(104, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(104, 9).Value = 100103002
$ws.Cells.Item(104, 10).Value = "Ciruela"
$ws.Cells.Item(104, 11).Value = "Black Cat"
$ws.Cells.Item(104, 12).Value = "Primera"
$ws.Cells.Item(104, 13).Value = 280
$ws.Cells.Item(104, 14).Value = 7500
$ws.Cells.Item(104, 15).Value = 7500
$ws.Cells.Item(104, 16).Value = 7500
$ws.Cells.Item(104, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(104, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(104, 19).Value = 500
$ws.Cells.Item(104, 20).Value = 15
